# Update "Pais" COVID-19 stats sheet with newer figures.
# - Bumps the "Datos actualizados..." timestamp in A1.
# - Refreshes Casos totales / Nuevos casos / Casos activos / Recuperados /
#   Casos criticos / Muertes hoy / Muertes for the countries whose counts
#   moved, which also re-orders a few neighbouring rows by rank
#   (Dinamarca now outranks Chile; Polonia now outranks Malasia; Oman now
#   outranks San Marino, Cuba and Vietnam).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 3 de Abril de 2020 a las 10:20"

function Set-Row($row, $pais, $casosTotales, $nuevosCasos, $casosActivos, $recuperados, $casosCriticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($row, 1).Value = $pais
    $ws.Cells.Item($row, 2).Value = $casosTotales
    $ws.Cells.Item($row, 3).Value = $nuevosCasos
    $ws.Cells.Item($row, 4).Value = $casosActivos
    $ws.Cells.Item($row, 5).Value = $recuperados
    $ws.Cells.Item($row, 6).Value = $casosCriticos
    $ws.Cells.Item($row, 7).Value = $muertesHoy
    $ws.Cells.Item($row, 8).Value = $muertes
}

# Row 17 - Austria
Set-Row 17 "Austria" 11224 95 1749 9317 227 0 158

# Rows 28-32 - Dinamarca/Chile swap rank (Ecuador stays at 30), Polonia/Malasia swap rank
Set-Row 28 "Dinamarca" 3672 286 1089 2460 153 0 123
Set-Row 29 "Chile" 3404 0 335 3051 31 0 18
Set-Row 31 "Polonia" 3149 203 56 3034 50 2 59
Set-Row 32 "Malasia" 3116 0 767 2299 105 0 50

# Row 58 - Ucrania
Set-Row 58 "Ucrania" 942 45 19 900 16 1 23

# Row 77 - Kazajistan
Set-Row 77 "Kazajistan" 453 18 29 421 6 0 3

# Rows 94-97 - Oman jumps ahead of San Marino, Cuba and Vietnam
Set-Row 94 "Oman" 252 21 57 194 3 0 1
Set-Row 95 "San Marino" 245 0 21 194 15 0 30
Set-Row 96 "Cuba" 233 0 13 214 7 0 6
Set-Row 97 "Vietnam" 233 0 85 148 3 0 0
